$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 0
$ws1.Range("F3").Value = 646
$ws1.Range("F4").Value = 351
$ws1.Range("F5").Value = 5030
$ws1.Range("F6").Value = 528
$ws1.Range("F8").Value = 242
$ws1.Range("F12").Value = 0

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14
$ws2.Range("F4").Value = 9
$ws2.Range("F6").Value = 3

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 646
$ws4.Range("F4").Value = 0
$ws4.Range("F7").Value = 5030
$ws4.Range("F8").Value = 528
$ws4.Range("F9").Value = 9
$ws4.Range("F10").Value = 9484
$ws4.Range("F11").Value = 242
$ws4.Range("F12").Value = 530
$ws4.Range("F14").Value = 7
$ws4.Range("F15").Value = 3
$ws4.Range("F16").Value = 673
$ws4.Range("F18").Value = 74
